# Fruta / hortaliza, semanal
# Inserts two new weekly records (Apio, Americana (o) - Primera/Segunda) for
# Terminal La Palmera de La Serena, pushing the existing data block down by
# two rows (old row 328 -> new row 330, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows right before the current row 328, shifting
# everything from row 328 onward down by two rows.
$ws.Rows("328:329").Insert()

# ---- New row 328: Apio, Americana (o), Primera ----
$ws.Range("A328").Value = 8
$ws.Range("B328").Value = "Terminal La Palmera de La Serena"
$ws.Range("C328").Value = "Coquimbo"
$ws.Range("D328").Value = 44736
$ws.Range("E328").Value = 4
$ws.Range("F328").Value = 100112017
$ws.Range("G328").Value = "Apio"
$ws.Range("H328").Value = "Americana (o)"
$ws.Range("I328").Value = "Primera"
$ws.Range("J328").Value = 2500
$ws.Range("K328").Value = 8000
$ws.Range("L328").Value = 9000
$ws.Range("M328").Value = 8500
$ws.Range("N328").Value = "`$/docena de matas"
$ws.Range("O328").Value = "Provincia del Elquí"
$ws.Range("P328").Value = 1417
$ws.Range("Q328").Value = 6
$ws.Range("R328").Value = "Hortaliza"

# ---- New row 329: Apio, Americana (o), Segunda ----
$ws.Range("A329").Value = 8
$ws.Range("B329").Value = "Terminal La Palmera de La Serena"
$ws.Range("C329").Value = "Coquimbo"
$ws.Range("D329").Value = 44736
$ws.Range("E329").Value = 4
$ws.Range("F329").Value = 100112017
$ws.Range("G329").Value = "Apio"
$ws.Range("H329").Value = "Americana (o)"
$ws.Range("I329").Value = "Segunda"
$ws.Range("J329").Value = 1400
$ws.Range("K329").Value = 6500
$ws.Range("L329").Value = 7000
$ws.Range("M329").Value = 6750
$ws.Range("N329").Value = "`$/docena de matas"
$ws.Range("O329").Value = "Provincia del Elquí"
$ws.Range("P329").Value = 1125
$ws.Range("Q329").Value = 6
$ws.Range("R329").Value = "Hortaliza"
